$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update attendance tracking values: mark the specific
# Real/Invalid/Absent indicator cells as 1 (from 0) for the
# corresponding attendance dates.

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("H4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
